$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1653846153846154
$ws.Range("C2").Value = 0.6230769230769231
$ws.Range("J2").Value = 0.003846153846153846
$ws.Range("P2").Value = 0.1153846153846154
$ws.Range("S2").Value = 0.09230769230769231
$ws.Range("B3").Value = 0.01183431952662722
$ws.Range("C3").Value = 0.04733727810650887
$ws.Range("J3").Value = 0.01183431952662722
$ws.Range("P3").Value = 0.6923076923076923
$ws.Range("S3").Value = 0.2366863905325444
$ws.Range("J4").Value = 0.03389830508474576
$ws.Range("O4").Value = 0.01694915254237288
$ws.Range("P4").Value = 0.6949152542372882
$ws.Range("S4").Value = 0.2542372881355932
$ws.Range("B6").Value = 0.06930693069306931
$ws.Range("D6").Value = 0.02475247524752475
$ws.Range("F6").Value = 0.08415841584158416
$ws.Range("J6").Value = 0.2722772277227723
$ws.Range("O6").Value = 0.009900990099009901
$ws.Range("Q6").Value = 0.202970297029703
$ws.Range("R6").Value = 0.03465346534653466
$ws.Range("S6").Value = 0.301980198019802
$ws.Range("B7").Value = 0.1348314606741573
$ws.Range("D7").Value = 0.0449438202247191
$ws.Range("F7").Value = 0.0449438202247191
$ws.Range("J7").Value = 0.06741573033707865
$ws.Range("O7").Value = 0.02808988764044944
$ws.Range("Q7").Value = 0.1910112359550562
$ws.Range("R7").Value = 0.06741573033707865
$ws.Range("S7").Value = 0.4213483146067415
$ws.Range("B8").Value = 0.09631728045325778
$ws.Range("D8").Value = 0.0198300283286119
$ws.Range("E8").Value = 0.0028328611898017
$ws.Range("F8").Value = 0.07082152974504249
$ws.Range("J8").Value = 0.1019830028328612
$ws.Range("O8").Value = 0.03682719546742209
$ws.Range("Q8").Value = 0.2436260623229462
$ws.Range("R8").Value = 0.06515580736543909
$ws.Range("S8").Value = 0.3626062322946176
$ws.Range("B9").Value = 0.106145251396648
$ws.Range("D9").Value = 0.01675977653631285
$ws.Range("F9").Value = 0.03910614525139665
$ws.Range("J9").Value = 0.0782122905027933
$ws.Range("O9").Value = 0.0223463687150838
$ws.Range("Q9").Value = 0.2290502793296089
$ws.Range("R9").Value = 0.05586592178770949
$ws.Range("S9").Value = 0.4525139664804469
$ws.Range("B10").Value = 0.1251212415130941
$ws.Range("D10").Value = 0.03685741998060136
$ws.Range("E10").Value = 0.0009699321047526673
$ws.Range("F10").Value = 0.07177497575169738
$ws.Range("J10").Value = 0.1086323957322987
$ws.Range("O10").Value = 0.02715809893307469
$ws.Range("Q10").Value = 0.2172647914645975
$ws.Range("R10").Value = 0.0504364694471387
$ws.Range("S10").Value = 0.3617846750727449
$ws.Range("G11").Value = 0.1454545454545454
$ws.Range("J11").Value = 0.07636363636363637
$ws.Range("K11").Value = 0.1890909090909091
$ws.Range("L11").Value = 0.5781818181818181
$ws.Range("S11").Value = 0.01090909090909091
$ws.Range("G12").Value = 0.7378048780487805
$ws.Range("J12").Value = 0.2134146341463415
$ws.Range("K12").Value = 0.006097560975609756
$ws.Range("L12").Value = 0.01219512195121951
$ws.Range("S12").Value = 0.03048780487804878
$ws.Range("G13").Value = 0.5681818181818182
$ws.Range("J13").Value = 0.3863636363636364
$ws.Range("S13").Value = 0.04545454545454546
$ws.Range("F15").Value = 0.0187793427230047
$ws.Range("H15").Value = 0.1032863849765258
$ws.Range("I15").Value = 0.07981220657276995
$ws.Range("J15").Value = 0.3427230046948357
$ws.Range("K15").Value = 0.07981220657276995
$ws.Range("M15").Value = 0.0187793427230047
$ws.Range("O15").Value = 0.1032863849765258
$ws.Range("S15").Value = 0.2535211267605634
$ws.Range("F16").Value = 0.02717391304347826
$ws.Range("H16").Value = 0.1521739130434783
$ws.Range("I16").Value = 0.09782608695652174
$ws.Range("J16").Value = 0.4076086956521739
$ws.Range("K16").Value = 0.1413043478260869
$ws.Range("M16").Value = 0.02173913043478261
$ws.Range("O16").Value = 0.05434782608695652
$ws.Range("S16").Value = 0.09782608695652174
$ws.Range("F17").Value = 0.03080568720379147
$ws.Range("H17").Value = 0.1729857819905213
$ws.Range("I17").Value = 0.1018957345971564
$ws.Range("J17").Value = 0.3862559241706161
$ws.Range("K17").Value = 0.1303317535545024
$ws.Range("M17").Value = 0.02606635071090047
$ws.Range("O17").Value = 0.06398104265402843
$ws.Range("S17").Value = 0.08767772511848342
$ws.Range("H18").Value = 0.09615384615384616
$ws.Range("I18").Value = 0.09615384615384616
$ws.Range("J18").Value = 0.4423076923076923
$ws.Range("K18").Value = 0.1057692307692308
$ws.Range("M18").Value = 0.02884615384615385
$ws.Range("O18").Value = 0.08653846153846154
$ws.Range("S18").Value = 0.1442307692307692
$ws.Range("F19").Value = 0.01945525291828794
$ws.Range("H19").Value = 0.2159533073929961
$ws.Range("I19").Value = 0.08657587548638132
$ws.Range("J19").Value = 0.3706225680933852
$ws.Range("K19").Value = 0.1079766536964981
$ws.Range("M19").Value = 0.02140077821011673
$ws.Range("N19").Value = 0.0009727626459143969
$ws.Range("O19").Value = 0.06906614785992218
$ws.Range("S19").Value = 0.1079766536964981
